$d = $word.ActiveDocument

function Get-ParagraphByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Remove the "No sistema atual..." paragraph, the two blank paragraphs
#    after it, and the "Com a implantacao do sistema, havera mudancas:"
#    paragraph that used to introduce the bullet list (paragraphs 4-7).
# ---------------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(4)
$pEnd   = $d.Paragraphs.Item(7)
$d.Range($pStart.Range.Start, $pEnd.Range.End).Delete()

# ---------------------------------------------------------------------------
# 2) "Atualmente sao aplicadas as seguintes: " -> add "formas de trabalho
#    pela oficina" before the colon, then add a new blank paragraph,
#    indented (ind left=720), right after it.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Atualmente são aplicadas as seguintes: ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Atualmente são aplicadas as seguintes formas de trabalho pela oficina: ", 2) | Out-Null

$pAtualmente = Get-ParagraphByText $d "Atualmente são aplicadas as seguintes formas de trabalho pela oficina:"
$pAtualmente.Range.InsertParagraphAfter() | Out-Null
$pIndent1 = $pAtualmente.Next()
$pIndent1.Format.LeftIndent = 36

# ---------------------------------------------------------------------------
# 3) "Proposta incluindo o Software:" -> full rewrite describing the new
#    changes that will happen once the software is implemented, with a
#    _GoBack bookmark right before the trailing ": ", then add a new
#    blank paragraph, indented (ind left=720), right after it.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Proposta incluindo o Software:", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Com a implantação do Software System Control Mecânica ocorrerá novas mudanças na forma de trabalho da oficina, tais como: ", 2) | Out-Null

$pComA = Get-ParagraphByText $d "ocorrerá novas mudanças na forma de trabalho da oficina, tais como:"
$bmPos = $pComA.Range.End - 3
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$pComA.Range.InsertParagraphAfter() | Out-Null
$pIndent2 = $pComA.Next()
$pIndent2.Format.LeftIndent = 36

# ---------------------------------------------------------------------------
# 4) "Organização e acessibilidade " -> "Organização e acessibilidade. "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Organização e acessibilidade ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Organização e acessibilidade. ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "...(entrada, pausa e saída) " -> "...(entrada, pausa e saída)." and
#    remove the following blank, indented (ind left=1440) paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Controle de funcionários referente ao ponto (entrada, pausa e saída) ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Controle de funcionários referente ao ponto (entrada, pausa e saída).", 2) | Out-Null

$pSaida = Get-ParagraphByText $d "Controle de funcionários referente ao ponto (entrada, pausa e saída)."
$pBlank1440 = $pSaida.Next()
$d.Range($pSaida.Range.End, $pBlank1440.Range.End).Delete()

# ---------------------------------------------------------------------------
# 6) Table: merge "João Paulo" + " " + "Souza" (with the old _GoBack
#    bookmark in between) into a single "João Paulo" + " Souza" run. The
#    _GoBack bookmark previously anchored here is implicitly dropped since
#    it moved to the "tais como" paragraph above (step 3).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "João Paulo Souza", $true, $false, $false, $false, $false,
    $true, 1, $false, "João Paulo Souza", 2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Theme colors: swap dk1/lt1 sysClr lastClr values.
# ---------------------------------------------------------------------------
$theme = $d.DocumentTheme
if ($theme -ne $null) {
    $scheme = $theme.ThemeColorScheme
}
